$d = $word.ActiveDocument

# Find the paragraph containing "LOB1232: Licenciamento Ambiental (Requisito)"
# and remove the following three paragraphs:
#   1. the blank paragraph right after it
#   2. "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3. "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

$target1 = "Ver no Jupiter Salvar em pdf Salvar em docx"
$target2 = "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

$paras = $d.Paragraphs
$count = $paras.Count

for ($i = $count; $i -ge 1; $i--) {
    $p = $paras.Item($i)
    $text = $p.Range.Text.Trim()
    if ($text -eq $target1 -or $text -eq $target2) {
        $p.Range.Delete()
    }
}

# Now remove the blank paragraph that immediately follows "LOB1232: Licenciamento Ambiental (Requisito)"
$paras = $d.Paragraphs
$count = $paras.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $paras.Item($i)
    $text = $p.Range.Text.Trim()
    if ($text -eq "LOB1232: Licenciamento Ambiental (Requisito)") {
        $next = $paras.Item($i + 1)
        if ($next.Range.Text.Trim() -eq "") {
            $next.Range.Delete()
        }
        break
    }
}
